$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.522.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.697.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.598"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.112.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.688.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.934"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.461.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.68%  "
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "281.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "31.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("E30").Value = "  -4.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0848"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.125"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.57%  "
$ws.Range("E43").Value = "  -4.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.119.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "93.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.947.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("E51").Value = "  -2.68%  "
